$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.065.75"
$ws.Range("E2").Value = "  -3.92%  "

$ws.Range("D3").Value = "2.332.78"
$ws.Range("E3").Value = "  -5.68%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'307.37"
$ws.Range("E5").Value = "  -4.01%  "

$ws.Range("D6").Value = "'84.97"
$ws.Range("E6").Value = "  -7.73%  "

$ws.Range("D7").Value = "'0.529"
$ws.Range("E7").Value = "  -3.69%  "

$ws.Range("D9").Value = "'0.483"
$ws.Range("E9").Value = "  -4.93%  "

$ws.Range("D10").Value = "'0.0816"
$ws.Range("E10").Value = "  -4.35%  "

$ws.Range("D11").Value = "'30.09"

$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").Value = "2.692.68"
$ws.Range("E13").Value = "  -5.68%  "

$ws.Range("D14").Value = "'6.40"
$ws.Range("E14").Value = "  -6.81%  "

$ws.Range("D15").Value = "'14.69"
$ws.Range("E15").Value = "  -5.02%  "

$ws.Range("D16").Value = "2.325.80"
$ws.Range("E16").Value = "  -5.86%  "

$ws.Range("D17").Value = "'0.752"
$ws.Range("E17").Value = "  -4.76%  "

$ws.Range("D18").Value = "40.038.29"
$ws.Range("E18").Value = "  -3.84%  "

$ws.Range("D19").Value = "0.0₃0903"
$ws.Range("E19").Value = "  -3.75%  "

$ws.Range("D20").Value = "'6.08"
$ws.Range("E20").Value = "  -5.53%  "

$ws.Range("D21").Value = "'67.57"
$ws.Range("E21").Value = "  -5.31%  "

$ws.Range("D22").Value = "'10.66"
$ws.Range("E22").Value = "  -5.01%  "

$ws.Range("D23").Value = "'235.43"
$ws.Range("E23").Value = "  -1.53%  "

$ws.Range("D24").Value = "'2.56"
$ws.Range("E24").Value = "  -6.90%  "

$ws.Range("E26").Value = "  -7.06%  "

$ws.Range("D27").Value = "'23.36"
$ws.Range("E27").Value = "  -5.97%  "

$ws.Range("D28").Value = "'2.21"
$ws.Range("E28").Value = "  -1.13%  "

$ws.Range("D29").Value = "'9.27"
$ws.Range("E29").Value = "  -4.74%  "

$ws.Range("D30").Value = "'35.44"
$ws.Range("E30").Value = "  -2.79%  "

$ws.Range("D31").Value = "'151.94"
$ws.Range("E31").Value = "  -2.34%  "

$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("D33").Value = "'5.11"
$ws.Range("E33").Value = "  -5.91%  "

$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.44"
$ws.Range("E34").Value = "  -4.47%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0724"
$ws.Range("E35").Value = "  -5.30%  "

$ws.Range("E36").Value = "  -1.95%  "

$ws.Range("D37").Value = "'0.0997"
$ws.Range("E37").Value = "  -3.30%  "

$ws.Range("E38").Value = "  -4.45%  "

$ws.Range("D39").Value = "'15.75"
$ws.Range("E39").Value = "  -7.86%  "

$ws.Range("D40").Value = "'1.70"
$ws.Range("E40").Value = "  -6.95%  "

$ws.Range("E41").Value = "  -4.44%  "

$ws.Range("D42").Value = "'2.26"
$ws.Range("E42").Value = "  -6.48%  "

$ws.Range("D43").Value = "1.939.40"
$ws.Range("E43").Value = "  -3.12%  "

$ws.Range("D44").Value = "'0.0267"
$ws.Range("E44").Value = "  -5.44%  "

$ws.Range("D45").Value = "'17.57"
$ws.Range("E45").Value = "  -5.77%  "

$ws.Range("D46").Value = "'9.26"
$ws.Range("E46").Value = "  -1.75%  "

$ws.Range("D47").Value = "'2.68"
$ws.Range("E47").Value = "  -9.27%  "

$ws.Range("D48").Value = "2.558.43"
$ws.Range("E48").Value = "  -6.32%  "

$ws.Range("D49").Value = "'92.89"
$ws.Range("E49").Value = "  -4.49%  "

$ws.Range("D50").Value = "'71.49"
$ws.Range("E50").Value = "  -5.44%  "

$ws.Range("D51").Value = "'50.62"
$ws.Range("E51").Value = "  -2.79%  "
